{"js": "// Split the existing \"FirstParagraph\" reference paragraph into two:\n//   1) a new FirstParagraph paragraph with the sentence \"Some random text.\"\n//   2) the original paragraph (unchanged runs), now styled as BodyText.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style,items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(\"A reference\") === 0) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find target paragraph ('A reference\u2026')\");\n}\n\n// Insert a new paragraph right before the target; it inherits the target's\n// style (FirstParagraph) automatically, then gets its own text.\nconst inserted = target.insertParagraph(\"Some random text.\", \"Before\");\ninserted.style = \"FirstParagraph\";\n\n// The original paragraph (with its original runs untouched) becomes BodyText.\ntarget.style = \"BodyText\";\n\nawait context.sync();\n", "ps1": "# Split the existing \"FirstParagraph\" reference paragraph into two:\n#   1) a new FirstParagraph paragraph with the sentence \"Some random text.\"\n#   2) the original paragraph (unchanged runs), now styled as BodyText.\n$d = $word.ActiveDocument\n\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.StartsWith(\"A reference\")) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find target paragraph ('A reference\u2026')\"\n}\n\n$target = $d.Paragraphs.Item($targetIndex)\n\n# Inserts an empty paragraph immediately before $target, inheriting its\n# (FirstParagraph) style; this shifts $target to $targetIndex + 1.\n$target.Range.InsertParagraphBefore()\n\n$newPara = $d.Paragraphs.Item($targetIndex)\n$newPara.Range.Text = \"Some random text.\"\n\n$target = $d.Paragraphs.Item($targetIndex + 1)\n$target.Style = \"Body Text\"\n"}
